$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (Id=1): Supplier/DRX/I import record
$ws.Range("B2").Value = "SupplierDumb"
$ws.Range("F2").Value = "Feb 18, 2022 (11:16:32 EST)"

# Row 3 (Id=2): Plant/DRX/I import record
$ws.Range("C3").Value = "DRXX"
$ws.Range("F3").Value = "Feb 18, 2022 (11:14:55 EST)"

# Row 4 (Id=3): Solicitation/DRX/I import record
$ws.Range("F4").Value = "Feb 18, 2022 (11:16:05 EST)"
